# Append a new "mudancas" service, the quote created for it, and the
# corresponding quote_item line, mirroring the existing rows already present
# in each sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. services!A21:N21 - new service "mudancas"
# ---------------------------------------------------------------------
$svc = $wb.Worksheets.Item("services")
$svc.Range("A21").Value = 20
$svc.Range("B21").Value = "mudancas"
$svc.Range("D21").Value = 52.13
$svc.Range("G21").Value = "unidade"
$svc.Range("M21").Value = "2025-09-25T14:04:10.329003"
$svc.Range("N21").Value = "2025-09-25T14:04:10.329003"

# ---------------------------------------------------------------------
# 2. quotes!A15:S15 - new quote for the "mudancas" service
# ---------------------------------------------------------------------
$quo = $wb.Worksheets.Item("quotes")
$quo.Range("A15").Value = 15
$quo.Range("B15").Value = "ORC202509013"
$quo.Range("C15").Value = 1
$quo.Range("D15").Value = "Orçamento - mudancas"
$quo.Range("E15").Value = "gerais"
$quo.Range("H15").Value = "pendente"
$quo.Range("M15").Value = 52.13
$quo.Range("R15").Value = "2025-09-25T14:04:10.698587"
$quo.Range("S15").Value = "2025-09-25T14:04:10.698587"

# ---------------------------------------------------------------------
# 3. quote_items!A16:O16 - line item linking the quote and the service
# ---------------------------------------------------------------------
$qi = $wb.Worksheets.Item("quote_items")
$qi.Range("A16").Value = 17
$qi.Range("B16").Value = 15
$qi.Range("C16").Value = 20
$qi.Range("D16").Value = 1
$qi.Range("E16").Value = 52.13
$qi.Range("F16").Value = ""
$qi.Range("G16").Value = 52.13
$qi.Range("H16").Value = "mudancas"
$qi.Range("I16").Value = ""
$qi.Range("J16").Value = "unidade"
$qi.Range("K16").Value = ""
$qi.Range("L16").Value = ""
$qi.Range("M16").Value = ""
$qi.Range("N16").Value = ""
$qi.Range("O16").Value = "2025-09-25T14:04:10.698587"
